# Updated cryptos list (price / 1h volume refresh, plus two rank swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings are not
# reinterpreted as floating point numbers, then restore the original "Normal"
# cell style (no explicit format) once the text values have been written.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.210.83'
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '2.550.66'
$ws.Range("E3").Value = '  +0.61%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '318.51'
$ws.Range("E5").Value = '  +4.26%  '

$ws.Range("D6").Value = '96.30'
$ws.Range("E6").Value = '  -2.65%  '

$ws.Range("E7").Value = '  -0.70%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("E9").Value = '  -2.13%  '

$ws.Range("D10").Value = '36.69'
$ws.Range("E10").Value = '  -0.98%  '

$ws.Range("D11").Value = '0.0815'
$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("E12").Value = '  -1.08%  '

$ws.Range("E13").Value = '  +0.35%  '

$ws.Range("D14").Value = '2.940.97'
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("D15").Value = '15.68'
$ws.Range("E15").Value = '  +3.72%  '

$ws.Range("D16").Value = '2.526.62'
$ws.Range("E16").Value = '  -1.63%  '

$ws.Range("D17").Value = '0.856'
$ws.Range("E17").Value = '  -1.52%  '

$ws.Range("D18").Value = '43.091.32'
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").Value = '13.10'
$ws.Range("E19").Value = '  +0.50%  '

$ws.Range("D20").Value = '6.65'
$ws.Range("E20").Value = '  +2.44%  '

$ws.Range("D21").Value = '0.0₃0973'
$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").Value = '70.62'
$ws.Range("E22").Value = '  -1.24%  '

$ws.Range("D23").Value = '252.89'
$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").Value = '2.98'
$ws.Range("E24").Value = '  +1.67%  '

$ws.Range("D25").Value = '2.03'
$ws.Range("E25").Value = '  -0.90%  '

$ws.Range("D26").Value = '27.13'
$ws.Range("E26").Value = '  +0.87%  '

$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("E28").Value = '  +3.72%  '

$ws.Range("D29").Value = '40.01'
$ws.Range("E29").Value = '  +4.09%  '

$ws.Range("D30").Value = '10.26'
$ws.Range("E30").Value = '  -1.79%  '

$ws.Range("D31").Value = '6.12'
$ws.Range("E31").Value = '  +0.20%  '

$ws.Range("D32").Value = '155.52'
$ws.Range("E32").Value = '  -1.75%  '

$ws.Range("E33").Value = '  +1.92%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '3.35'
$ws.Range("E34").Value = '  +0.78%  '

$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").Value = '19.10'
$ws.Range("E35").Value = '  +4.51%  '

$ws.Range("D36").Value = '0.0794'
$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("D37").Value = '2.63'
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("E38").Value = '  -2.40%  '

$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").Value = '24.15'
$ws.Range("E39").Value = '  -0.24%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.119'
$ws.Range("E40").Value = '  -0.80%  '

$ws.Range("D41").Value = '2.29'
$ws.Range("E41").Value = '  +9.76%  '

$ws.Range("D42").Value = '3.84'
$ws.Range("E42").Value = '  -1.80%  '

$ws.Range("D43").Value = '3.36'
$ws.Range("E43").Value = '  -2.42%  '

$ws.Range("D44").Value = '0.0305'
$ws.Range("E44").Value = '  +0.66%  '

$ws.Range("E45").Value = '  +0.26%  '

$ws.Range("D46").Value = '2.024.48'
$ws.Range("E46").Value = '  -0.97%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("D48").Value = '8.88'
$ws.Range("E48").Value = '  -0.91%  '

$ws.Range("D49").Value = '2.790.15'
$ws.Range("E49").Value = '  +0.43%  '

$ws.Range("D50").Value = '74.85'
$ws.Range("E50").Value = '  +2.46%  '

$ws.Range("D51").Value = '103.21'
$ws.Range("E51").Value = '  +0.16%  '

$ws.Range("D2:D51").Style = "Normal"